# Refresh the crypto price/volume table (cryptos.xlsx) with the latest
# scrape. Price text in column D and the volume-change percentages in
# column E are stored as plain text (not numbers), matching the site's
# formatting (thousands-dot notation, padded "  +x.xx%  " strings).
#
# Some column-D prices (e.g. "602.80", "6.94") look numeric, so a plain
# Range.Value assignment would get auto-coerced into a float and lose the
# literal text (trailing zeros, exact precision). To keep them as text -
# exactly like the original cells - we briefly force a text number format,
# assign the value, then clear the formatting again so the cell's style
# index stays untouched (identical to how it started).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.399.09'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.497.24'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -4.82%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.80'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.35'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.496.42'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -4.79%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -3.24%  '
$ws.Range("E10").Value = '  -4.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.94'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.97%  '
$ws.Range("E12").Value = '  -4.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000217'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.089.04'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.50'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.501.99'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.304.85'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.15%  '
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.97'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -5.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '447.32'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.98'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -12.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.618'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -5.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.36'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000128'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.638.30'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.09'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -9.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.19'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.22%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.53'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.44%  '
$ws.Range("E33").Value = '  +1.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.68'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.07'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.488.30'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.29%  '
$ws.Range("E37").Value = '  -6.95%  '
$ws.Range("E38").Value = '  -3.76%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '174.97'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("E42").Value = '  -1.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0873'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("E45").Value = '  -4.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.43'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.41'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.35%  '
$ws.Range("E48").Value = '  +6.03%  '
$ws.Range("E49").Value = '  -5.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.52'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.995'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.08%  '
